# Generate Report for Handoff
# Updates Priority ("low" -> "ht") and Latest Handoff Datetime for the rows
# that were just handed off, on both the zh-cn and de-de localization-status
# worksheets, and refreshes the Overview sheet's Latest HO Xliff Generate
# Date to match the newest handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-13 10:37:39"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-13 10:37:46"
}

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" tracks the newest handoff date across
# locales, which is the de-de one that was just regenerated.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2016-08-13 10:37:46"
}
